$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "dog" sheet — append a new field-test result as row 5
# ---------------------------------------------------------------------
$dog = $wb.Worksheets.Item("dog")

# Clone the formatting of the row above (keeps date / time number formats
# and reuses the existing style records instead of minting new ones).
$dog.Range("A4:L4").Copy()
$dog.Range("A5:L5").PasteSpecial(-4122)

$dog.Cells.Item(5, 1).Value = 45779                                   # Date -> 2025-05-02
$dog.Cells.Item(5, 2).Value = "PRESENCE"                              # Type
$dog.Cells.Item(5, 3).Value = 0.61458333333333337                     # Time target placed -> 14:45
$dog.Cells.Item(5, 4).Value = 0.70833333333333337                     # Time of search -> 17:00
$dog.Cells.Item(5, 5).Value = 14                                      # Temperature degrees
$dog.Cells.Item(5, 6).Value = 0                                       # Wind kmh
$dog.Cells.Item(5, 7).Value = "Sunny, cool"                           # Conditions
$dog.Cells.Item(5, 8).Value = $true                                   # Found
$dog.Cells.Item(5, 9).Value = "19 minutes 51 seconds"                 # Search time mins
$dog.Cells.Item(5, 10).Value = 1191                                   # Search time s
$dog.Cells.Item(5, 11).Value = "Primary sweeps"                       # Search stage
$dog.Cells.Item(5, 12).Value = "Worked downhill, transmitter near end of search. Really clear alert even when I tried to pull her off target and good focus with duration."  # Notes

# ---------------------------------------------------------------------
# 2. "human" sheet — drop the stale blank placeholder rows (4:34) and
#    clear the leftover number formatting on F3
# ---------------------------------------------------------------------
$human = $wb.Worksheets.Item("human")

$human.Range("A4:H34").EntireRow.Delete()
$human.Range("F3").ClearFormats()

# ---------------------------------------------------------------------
# 3. Selection / active-tab bookkeeping
# ---------------------------------------------------------------------
$human.Range("H7").Select()

$dog.Activate()
$dog.Range("A6").Select()
